$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 27 (the SCYOCAN / 2-RAP entry). This shifts rows 28-30 up to
# 27-29, yielding the target data and shrinking the used range to A1:K29.
$ws.Rows("27").Delete()
